$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.147174835205078
$ws.Range("B1").Value = 2.636829614639282
$ws.Range("C1").Value = 5.839153289794922
$ws.Range("D1").Value = 2.088707208633423
$ws.Range("E1").Value = 1.202476382255554
